$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> A팀)
$ws.Name = "A팀"

# Shorten the allowance/overtime header labels (drop the trailing unit in parentheses)
$ws.Range("E1").Value = "초과근무시간"
$ws.Range("F1").Value = "직책수당"
$ws.Range("G1").Value = "자격증수당"
$ws.Range("H1").Value = "위험수당"
$ws.Range("I1").Value = "수당합계"

# Re-fit the affected columns to the new, shorter header text
$ws.Columns.Item(5).ColumnWidth = 83/7
$ws.Columns.Item(6).ColumnWidth = 55/7
$ws.Columns.Item(7).ColumnWidth = 69/7
$ws.Columns.Item(8).ColumnWidth = 55/7
$ws.Columns.Item(9).ColumnWidth = 55/7

# Move the active selection to K5
[void]$ws.Range("K5").Select()
